$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the existing "charging_station1/charging_station" row
# (currently row 5), pushing it down to row 6, then populate the new row 5
# with the battery entry.
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "bat1"
$ws.Range("B5").Value = "bat"
